# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;   Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 24;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 26;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 40;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 48;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 63;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 64;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 76;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 89;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 105; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 106; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 107; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 114; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 116; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 121; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 125; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 151; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 153; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 162; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 163; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 168; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 169; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 173; Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 183; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 184; Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 197; Tag = "sd"; Act = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
